$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# The "task_title" calculate row (old row 24) and the "namee" note row
# (old row 26) are no longer used by the form -- drop them. Deleting row
# 24 first shifts the old row 26 up to row 25, so it is removed second.
$ws.Rows("24:24").Delete()
$ws.Rows("25:25").Delete()

# Bold the note labels (markdown-style **text**) and repoint the title
# note at the renamed ${my_field_title} field. Each note cell is rich
# text: run 1 is the bold label, run 2 is the plain-text value
# reference. Replacing .Text resets formatting, so re-apply Bold/Font
# after each edit.

# title note (now row 25): run1 = chars 1-32, run2 = chars 33-49
$titleCell = $ws.Range("C25")
$newTitleLabel = " **What is the title of the task?**:"
$newTitleValue = "  `${my_field_title} `n"
$titleCell.Characters(1, 32).Text = $newTitleLabel
$titleLabelRun = $titleCell.Characters(1, $newTitleLabel.Length)
$titleLabelRun.Font.Name = "Arial"
$titleLabelRun.Font.Bold = $true
$titleValueRun = $titleCell.Characters($newTitleLabel.Length + 1, $newTitleValue.Length)
$titleValueRun.Text = $newTitleValue
$titleValueRun.Font.Name = "Arial"
$titleValueRun.Font.Bold = $false

# when note (now row 27): run1 = chars 1-45, run2 = chars 46-57
$whenCell = $ws.Range("C27")
$newWhenLabel = " **When should a reminder for this task appear?**"
$whenCell.Characters(1, 45).Text = $newWhenLabel
$whenLabelRun = $whenCell.Characters(1, $newWhenLabel.Length)
$whenLabelRun.Font.Name = "Arial"
$whenLabelRun.Font.Bold = $true
$whenValueRun = $whenCell.Characters($newWhenLabel.Length + 1, 12)
$whenValueRun.Font.Name = "Arial"
$whenValueRun.Font.Bold = $false

# notes note (now row 28): run1 = chars 1-24, run2 = chars 25-43
$notesCell = $ws.Range("C28")
$newNotesLabel = " **Notes about this task:** "
$notesCell.Characters(1, 24).Text = $newNotesLabel
$notesLabelRun = $notesCell.Characters(1, $newNotesLabel.Length)
$notesLabelRun.Font.Name = "Arial"
$notesLabelRun.Font.Bold = $true
$notesValueRun = $notesCell.Characters($newNotesLabel.Length + 1, 19)
$notesValueRun.Font.Name = "Arial"
$notesValueRun.Font.Bold = $false
